# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the aggregated "全部类型" sheet to match freshly generated data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 99
$ws1.Range("F3").Value = 4080
$ws1.Range("F11").Value = 90
$ws1.Range("F13").Value = 1528
$ws1.Range("F15").Value = 2950

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 99
$ws4.Range("F3").Value = 4080
$ws4.Range("F17").Value = 1528
$ws4.Range("F19").Value = 2950
